$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the player_1 chat lines (rows 5-6) to the new lines
$ws.Range("B5").Value = "I can be your angle…. or yuor devil"
$ws.Range("B6").Value = "All your base are belong to us!"

# Insert 6 new rows after row 8 (after the "botch" category block) for the
# 3 new categories: fizzle, oncooldown, cooldownlistfull
$ws.Rows("9:14").Insert()

# Write the three new category names first (so shared strings are
# registered in the same order as the target workbook)
$ws.Range("B9").Value = "fizzle"
$ws.Range("B11").Value = "oncooldown"
$ws.Range("B13").Value = "cooldownlistfull"

# Then write the three new dialogue lines
$ws.Range("B10").Value = "Wow, I was WAY off…"
$ws.Range("B12").Value = "Oh no, it's on cooldown!"
$ws.Range("B14").Value = "My cooldown list is full!"

# fizzle category
$ws.Range("A9").Value = "NEW_CATEGORY"
$ws.Range("A10").Value = 100

# oncooldown category
$ws.Range("A11").Value = "NEW_CATEGORY"
$ws.Range("A12").Value = 100

# cooldownlistfull category
$ws.Range("A13").Value = "NEW_CATEGORY"
$ws.Range("A14").Value = 100

# Update the sheet view selection to match the target state
$ws.Range("B14").Select()

Write-Output "done"
